# Update cosinor analysis results for square_05 (CircaDB / CircadiPy re-run)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 ---
$ws.Range("E2").Value = [double]"22.34000000000005"
$ws.Range("H2").Value = [double]"1.708035422500241e-16"
$ws.Range("K2").Value = [double]"51.95602624330596"
$ws.Range("L2").Value = "[44.303730767852414, 59.6083217187595]"
$ws.Range("O2").Value = [double]"1.553500271144502"
$ws.Range("P2").Value = "[1.3899739268135018, 1.7170266154755023]"
$ws.Range("S2").Value = [double]"57.08333091349763"
$ws.Range("T2").Value = "[51.88371871980844, 62.282943107186824]"
$ws.Range("W2").Value = [double]"16.81649649649654"
$ws.Range("X2").Value = [double]"16.23507507507512"
$ws.Range("Y2").Value = [double]"17.39791791791796"

# --- Row 3 ---
$ws.Range("E3").Value = [double]"24.77000000000043"
$ws.Range("H3").Value = [double]"1.708035422500241e-16"
$ws.Range("K3").Value = [double]"52.25304121964751"
$ws.Range("L3").Value = "[45.24931122812851, 59.256771211166516]"
$ws.Range("O3").Value = [double]"-2.251631971942235"
$ws.Range("P3").Value = "[-2.377421467581467, -2.125842476303003]"
$ws.Range("S3").Value = [double]"55.71404856357878"
$ws.Range("T3").Value = "[52.160255664837486, 59.26784146232007]"
$ws.Range("W3").Value = [double]"8.876536536536696"
$ws.Range("X3").Value = [double]"8.380640640640786"
$ws.Range("Y3").Value = [double]"9.372432432432605"
